$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Uniformize style of E10/E11 (drop stray quote-prefix formatting) to match
# the rest of the MPN column (same look as E8/E9) before changing their text.
$ws.Range("E8").Copy()
$ws.Range("E10:E11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update MPN values (R1, R2, C1/C2, C3) to uniformize with other Harp boards
$ws.Range("E8").Value = "ERJ-2RKF1001X"
$ws.Range("E9").Value = "ERJ-2RKF1003X"
$ws.Range("E10").Value = "C1608X6S1C475K080AC"
$ws.Range("E11").Value = "GRM155R71E104KE14D"

# Update copyright notice
$ws.Range("B21").Value = "Copyright 2020-2023 Artur Silva and Filipe Carvalho"

# View: zoom + selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("B22").Select()

# Column F width
$ws.Columns("F").ColumnWidth = 18.29

Write-Host "done"
